$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New tracker rows for 2025-08-22 (auto-sync)
$ws.Range("A377:B391").NumberFormat = "@"

$ws.Range("A377").Value = "14427808"
$ws.Range("B377").Value = "2025-08-22"
$ws.Range("C377").Value = "Alycia Parks"
$ws.Range("D377").Value = "Diana Shnaider"
$ws.Range("E377").Value = "Gana Alycia Parks"
$ws.Range("F377").Value = 2.75

$ws.Range("A378").Value = "14427997"
$ws.Range("B378").Value = "2025-08-22"
$ws.Range("C378").Value = "Sorana Cirstea"
$ws.Range("D378").Value = "Anastasia Zakharova"
$ws.Range("E378").Value = "Gana Anastasia Zakharova"
$ws.Range("F378").Value = 3.4

$ws.Range("A379").Value = "14506182"
$ws.Range("B379").Value = "2025-08-22"
$ws.Range("C379").Value = "Arthur Cazaux"
$ws.Range("D379").Value = "Jan-Lennard Struff"
$ws.Range("E379").Value = "Gana Arthur Cazaux"
$ws.Range("F379").Value = 1.83

$ws.Range("A380").Value = "14506175"
$ws.Range("B380").Value = "2025-08-22"
$ws.Range("C380").Value = "Kimmer Coppejans"
$ws.Range("D380").Value = "Ignacio Buse"
$ws.Range("E380").Value = "Gana Kimmer Coppejans"
$ws.Range("F380").Value = 2.5

$ws.Range("A381").Value = "14506170"
$ws.Range("B381").Value = "2025-08-22"
$ws.Range("C381").Value = "Leandro Riedi"
$ws.Range("D381").Value = "Garrett Johns"
$ws.Range("E381").Value = "Gana Garrett Johns"
$ws.Range("F381").Value = 3.4

$ws.Range("A382").Value = "14506174"
$ws.Range("B382").Value = "2025-08-22"
$ws.Range("C382").Value = "Yuta Shimizu"
$ws.Range("D382").Value = "Martin Damm Jr"
$ws.Range("E382").Value = "Gana Yuta Shimizu"
$ws.Range("F382").Value = 2.38

$ws.Range("A383").Value = "14507069"
$ws.Range("B383").Value = "2025-08-22"
$ws.Range("C383").Value = "Marc-Andrea Huesler"
$ws.Range("D383").Value = "Zachary Svajda"
$ws.Range("E383").Value = "Gana Marc-Andrea Huesler"
$ws.Range("F383").Value = 2.75

$ws.Range("A384").Value = "14507980"
$ws.Range("B384").Value = "2025-08-22"
$ws.Range("C384").Value = "Jerome Kym"
$ws.Range("D384").Value = "Yibing Wu"
$ws.Range("E384").Value = "Gana Jerome Kym"
$ws.Range("F384").Value = 1.73

$ws.Range("A385").Value = "14506177"
$ws.Range("B385").Value = "2025-08-22"
$ws.Range("C385").Value = "Dalma Galfi"
$ws.Range("D385").Value = "Panna Udvardy"
$ws.Range("E385").Value = "Gana Panna Udvardy"
$ws.Range("F385").Value = 2.62

$ws.Range("A386").Value = "14506178"
$ws.Range("B386").Value = "2025-08-22"
$ws.Range("C386").Value = "Emerson Jones"
$ws.Range("D386").Value = "Victoria Jimenez Kasintseva"
$ws.Range("E386").Value = "Gana Victoria Jimenez Kasintseva"
$ws.Range("F386").Value = 2

$ws.Range("A387").Value = "14507067"
$ws.Range("B387").Value = "2025-08-22"
$ws.Range("C387").Value = "Alina Charaeva"
$ws.Range("D387").Value = "Xiyu Wang"
$ws.Range("E387").Value = "Gana Alina Charaeva"
$ws.Range("F387").Value = 3.25

$ws.Range("A388").Value = "14506181"
$ws.Range("B388").Value = "2025-08-22"
$ws.Range("C388").Value = "Tereza Valentova"
$ws.Range("D388").Value = "Arantxa Rus"
$ws.Range("E388").Value = "Gana Arantxa Rus"
$ws.Range("F388").Value = 3.4

$ws.Range("A389").Value = "14507066"
$ws.Range("B389").Value = "2025-08-22"
$ws.Range("C389").Value = "Varvara Gracheva"
$ws.Range("D389").Value = "Ena Shibahara"
$ws.Range("E389").Value = "Gana Ena Shibahara"
$ws.Range("F389").Value = 2.5

$ws.Range("A390").Value = "14507448"
$ws.Range("B390").Value = "2025-08-22"
$ws.Range("C390").Value = "Dominika Salkova"
$ws.Range("D390").Value = "Rebecca Marino"
$ws.Range("E390").Value = "Gana Dominika Salkova"
$ws.Range("F390").Value = 2

$ws.Range("A391").Value = "14503673"
$ws.Range("B391").Value = "2025-08-22"
$ws.Range("C391").Value = "Tom Paris"
$ws.Range("D391").Value = "Pedro Araujo"
$ws.Range("E391").Value = "Gana Pedro Araujo"
$ws.Range("F391").Value = 2.25

